# Update vm_pu.xlsx Case_1_91 results: B column slack-bus setpoint 1.05 -> 1.02 (380 kV case)
# and recomputed per-bus voltage magnitudes (columns C-F, I-N) for rows 2-25.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> column letter -> new value
$data = @{
    2 = @{ "B" = 1.02; "C" = 1.034966441954001; "D" = 1.057173074317487; "E" = 1.045865765346316; "F" = 1.061112509032066; "I" = 1.045160578730297; "J" = 1.040082795201068; "K" = 1.059908770096956; "L" = 1.048632805547242; "M" = 1.063837454410509; "N" = 1.041559832069417 }
    3 = @{ "B" = 1.02; "C" = 1.036047559176088; "D" = 1.057887188232904; "E" = 1.046786349117663; "F" = 1.06202428436622; "I" = 1.045353321118017; "J" = 1.040806599357579; "K" = 1.060436921429535; "L" = 1.04936463253108; "M" = 1.064563541913984; "N" = 1.042284664110852 }
    4 = @{ "B" = 1.02; "C" = 1.036747286081158; "D" = 1.058347638112024; "E" = 1.047382086634578; "F" = 1.062613325496037; "I" = 1.045475615375315; "J" = 1.04127456651464; "K" = 1.060776371286526; "L" = 1.049837603927576; "M" = 1.065031799895301; "N" = 1.042753295834953 }
    5 = @{ "B" = 1.02; "C" = 1.037041492583567; "D" = 1.058540819702333; "E" = 1.047632547997745; "F" = 1.062860732713861; "I" = 1.045526447036241; "J" = 1.041471208057541; "K" = 1.06091852433025; "L" = 1.050036304227525; "M" = 1.065228278948997; "N" = 1.0429502166314 }
    6 = @{ "B" = 1.02; "C" = 1.037090893577825; "D" = 1.058573232739654; "E" = 1.047674602324521; "F" = 1.062902260211291; "I" = 1.045534947811607; "J" = 1.041504219643683; "K" = 1.060942360065759; "L" = 1.050069658840645; "M" = 1.065261246534907; "N" = 1.042983275097781 }
    7 = @{ "B" = 1.02; "C" = 1.036751217120297; "D" = 1.058350220954221; "E" = 1.047385433260102; "F" = 1.06261663225145; "I" = 1.04547629687381; "J" = 1.041277194409834; "K" = 1.060778272911746; "L" = 1.049840259509215; "M" = 1.065034426739059; "N" = 1.04275592746206 }
    8 = @{ "B" = 1.02; "C" = 1.035331775131066; "D" = 1.057414749030952; "E" = 1.046176868500749; "F" = 1.061420841500314; "I" = 1.045226218145557; "J" = 1.040327487250354; "K" = 1.060087736698403; "L" = 1.048880247538046; "M" = 1.064083163436923; "N" = 1.04180487160948 }
    9 = @{ "B" = 1.02; "C" = 1.032831843166111; "D" = 1.055753908218551; "E" = 1.044047696218227; "F" = 1.0593065603271; "I" = 1.044767024134494; "J" = 1.038651063262285; "K" = 1.058853364707214; "L" = 1.047184242946442; "M" = 1.062394937657391; "N" = 1.040126066907058 }
    10 = @{ "B" = 1.02; "C" = 1.031166083284975; "D" = 1.054638409231926; "E" = 1.042628599470722; "F" = 1.057892289023352; "I" = 1.044448484011845; "J" = 1.03753149704367; "K" = 1.05801871051253; "L" = 1.046050680821997; "M" = 1.061261450250435; "N" = 1.039004910775998 }
    11 = @{ "B" = 1.02; "C" = 1.030444990325487; "D" = 1.054153439407956; "E" = 1.042014204282464; "F" = 1.057278777221265; "I" = 1.044307617433984; "J" = 1.037046251582086; "K" = 1.057654524495131; "L" = 1.045559154443852; "M" = 1.060768748753421; "N" = 1.0385189762102 }
    12 = @{ "B" = 1.02; "C" = 1.030177172607524; "D" = 1.053973007979798; "E" = 1.041786003412073; "F" = 1.057050723581238; "I" = 1.044254852722707; "J" = 1.036865939843967; "K" = 1.057518833566545; "L" = 1.045376476786665; "M" = 1.060585453560463; "N" = 1.038338408408731 }
    13 = @{ "B" = 1.02; "C" = 1.030234619144652; "D" = 1.054011724331234; "E" = 1.041834952685606; "F" = 1.057099649454703; "I" = 1.044266190860664; "J" = 1.036904620489722; "K" = 1.057547958542862; "L" = 1.045415666415464; "M" = 1.060624783839661; "N" = 1.03837714398544 }
    14 = @{ "B" = 1.02; "C" = 1.030422851854932; "D" = 1.054138530840935; "E" = 1.04199534086288; "F" = 1.057259929660373; "I" = 1.044303264871032; "J" = 1.037031348384762; "K" = 1.057643316722215; "L" = 1.045544056348866; "M" = 1.060753603303911; "N" = 1.038504051848627 }
    15 = @{ "B" = 1.02; "C" = 1.030538831975766; "D" = 1.054216621851712; "E" = 1.042094163027342; "F" = 1.057358661337916; "I" = 1.044326049027314; "J" = 1.037109420352679; "K" = 1.057702014943353; "L" = 1.045623147999946; "M" = 1.060832935623737; "N" = 1.038582234687693 }
    16 = @{ "B" = 1.02; "C" = 1.031213943850297; "D" = 1.054670554065548; "E" = 1.042669376638657; "F" = 1.057932982179216; "I" = 1.044457771043363; "J" = 1.037563691373453; "K" = 1.058042821947243; "L" = 1.046083287367964; "M" = 1.061294109376476; "N" = 1.039037150825422 }
    17 = @{ "B" = 1.02; "C" = 1.031637474733593; "D" = 1.054954772114247; "E" = 1.043030215451357; "F" = 1.058292938297566; "I" = 1.044539611000432; "J" = 1.037848518871201; "K" = 1.058255858825081; "L" = 1.046371737032035; "M" = 1.061582884868784; "N" = 1.039322382810887 }
    18 = @{ "B" = 1.02; "C" = 1.031884531700068; "D" = 1.055120363337009; "E" = 1.043240694616131; "F" = 1.058502786214585; "I" = 1.04458706348393; "J" = 1.038014608977677; "K" = 1.058379851800468; "L" = 1.046539918586866; "M" = 1.061751139946247; "N" = 1.039488708784366 }
    19 = @{ "B" = 1.02; "C" = 1.031968774906898; "D" = 1.055176793681286; "E" = 1.043312463915619; "F" = 1.058574320524406; "I" = 1.04460319547008; "J" = 1.038071233785849; "K" = 1.058422084743548; "L" = 1.046597252904824; "M" = 1.061808479557204; "N" = 1.03954541400626 }
    20 = @{ "B" = 1.02; "C" = 1.031592031961232; "D" = 1.054924297681652; "E" = 1.042991500032342; "F" = 1.058254329619558; "I" = 1.044530859654327; "J" = 1.037817964225739; "K" = 1.05823302968044; "L" = 1.046340795963855; "M" = 1.061551920892624; "N" = 1.039291784774324 }
    21 = @{ "B" = 1.02; "C" = 1.030367421236108; "D" = 1.054101197537222; "E" = 1.041948110184311; "F" = 1.057212735768086; "I" = 1.044292359656793; "J" = 1.036994032109563; "K" = 1.057615247592216; "L" = 1.045506251551317; "M" = 1.060715677015061; "N" = 1.038466682580036 }
    22 = @{ "B" = 1.02; "C" = 1.029597623411407; "D" = 1.053581992682888; "E" = 1.041292164147916; "F" = 1.056556872463546; "I" = 1.044139856396636; "J" = 1.036475588479913; "K" = 1.057224416911758; "L" = 1.044980945273772; "M" = 1.060188254923521; "N" = 1.037947502700967 }
    23 = @{ "B" = 1.02; "C" = 1.030005692483841; "D" = 1.053857392738006; "E" = 1.041639886361202; "F" = 1.056904649916329; "I" = 1.044220942665708; "J" = 1.036750463621228; "K" = 1.057431831463246; "L" = 1.045259476392976; "M" = 1.060468006827774; "N" = 1.038222768196511 }
    24 = @{ "B" = 1.02; "C" = 1.031612565536357; "D" = 1.054938068347954; "E" = 1.043008993838915; "F" = 1.058271775554083; "I" = 1.044534814886503; "J" = 1.037831770693574; "K" = 1.058243346015743; "L" = 1.046354777105322; "M" = 1.061565912745285; "N" = 1.039305610848927 }
    25 = @{ "B" = 1.02; "C" = 1.033477982411851; "D" = 1.056184738979002; "E" = 1.044598078726627; "F" = 1.059853993820909; "I" = 1.044887927379103; "J" = 1.039084803635248; "K" = 1.059174553364572; "L" = 1.047623212452533; "M" = 1.062832798737458; "N" = 1.040560423241103 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
